{"js": "// Insert a new warning paragraph about the \"Raw Report Tables\" file right\n// before the \"In all cases, the derived primary key...\" paragraph (i.e.\n// immediately after the \"Keys, Constraints, and Utility Functions:\" heading\n// and its bookmark end), matching the plain/default paragraph style used by\n// the surrounding body text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"In all cases, the derived primary key (UniqID) is a concatenation of the following fields: MembershipID, PatientDOB, PatientGenderCode, and FamilyMembershipID.\";\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not locate anchor paragraph for insertion.\");\n}\n\nconst warningText =\n  'The \"Raw Report Tables\" file has been left intact for testing purposes. Use the modular form in production, notably field ordering for gallbladder procedure tables differ between the two versions, and down-stream processing relies on the newer ordering within Modular_Report_Call.';\n\n// Insert a brand-new, plainly-styled paragraph immediately before the\n// anchor paragraph (so it lands right after the preceding heading).\nconst newParagraph = anchorParagraph.insertParagraph(warningText, \"Before\");\n\nawait context.sync();\n", "ps1": "# Insert a new warning paragraph about the \"Raw Report Tables\" file right\n# before the \"In all cases, the derived primary key...\" paragraph (i.e.\n# immediately after the \"Keys, Constraints, and Utility Functions:\" heading\n# and its bookmark end), matching the plain/default paragraph style used by\n# the surrounding body text.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"In all cases, the derived primary key (UniqID) is a concatenation of the following fields: MembershipID, PatientDOB, PatientGenderCode, and FamilyMembershipID.\"\n$warningText = 'The \"Raw Report Tables\" file has been left intact for testing purposes. Use the modular form in production, notably field ordering for gallbladder procedure tables differ between the two versions, and down-stream processing relies on the newer ordering within Modular_Report_Call.'\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Could not locate anchor paragraph for insertion.\"\n}\n\n# How many paragraphs precede the anchor paragraph tells us its (1-based)\n# position in $d.Paragraphs; the brand-new paragraph will land at that same\n# position once it is inserted immediately before the anchor.\n$newIndex = $d.Range(0, $rng.Start).Paragraphs.Count + 1\n\n# Insert a brand-new, plainly-styled paragraph immediately before the\n# anchor paragraph (so it lands right after the preceding heading).\n$rng.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs.Item($newIndex)\n$newPara.Range.Text = $warningText\n"}
